# LOB1046.xlsx content correction
#
# A new row is inserted (the "Docentes responsáveis" teacher-name value gets
# its own row instead of being crammed, mistakenly, into the "Objetivos"
# row), and several cells that had stale/duplicated text from earlier rows
# are corrected with their proper content (objectives, summarized syllabus,
# full syllabus, evaluation method text, recovery rule and bibliography).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above the current row 13 ("Programa resumido:" row),
#    shifting it (and everything below) down by one.
$ws.Rows.Item(13).Insert()

# The inserted row copies formatting (including custom row height) from the
# row that used to be #13; the target has no explicit height for this new
# row, so clear it back to the sheet default.
$ws.Rows.Item(13).AutoFit()

# 2) Row 10 ("Objetivos:") previously held the teacher's name by mistake in
#    B/C - replace it with the real objectives text.
$ws.Cells.Item(10, 2).Value = "Apresentar aos alunos os princípios fundamentais de engenharia do meio ambiente."
$ws.Cells.Item(10, 3).Value = "Apresentar aos alunos os princípios fundamentais de engenharia do meio ambiente."

# 3) New row 13 (under "Docentes responsáveis:") gets the teacher's name in
#    B/C, with no value in column A.
$ws.Cells.Item(13, 2).Value = "5840671 - Francisco José Moreira Chaves"
$ws.Cells.Item(13, 3).Value = "5840671 - Francisco José Moreira Chaves"

# 4) Row 14 ("Programa resumido:") previously held "Semestral" - fix to the
#    actual short syllabus text.
$ws.Cells.Item(14, 2).Value = "1 - Fundamentos da Engenharia e o Meio Ambiente. 2 - O meio ambiente aquático. 3 - O meio ambiente terrestre. 4 - O meio ambiente atmosférico ."
$ws.Cells.Item(14, 3).Value = "1 - Fundamentos da Engenharia e o Meio Ambiente. 2 - O meio ambiente aquático. 3 - O meio ambiente terrestre. 4 - O meio ambiente atmosférico ."

# 5) Row 16 ("Programa:") previously held a stray date - fix to the full
#    syllabus text.
$ws.Cells.Item(16, 2).Value = "1 - FUNDAMENTOS: A Engenharia e o Meio Ambiente; Os Ecossistemas. A crise energética. Fontes alternativas de energia. A sustentabilidade do meio ambiente. 2 - O MEIO AMBIENTE AQUÁTICO: Composição e Propriedades; Necessidade e Utilização; Requisitos de Qualidade; Poluição.3 - O MEIO AMBIENTE TERRESTRE: Composição e Propriedades; Necessidades e Utilização; Requisitos de Qualidade; Poluição.4 - O MEIO AMBIENTE ATMOSFÉRICO: Composição e Propriedades; Requisitos de Qualidade; Poluição."
$ws.Cells.Item(16, 3).Value = "1 - FUNDAMENTOS: A Engenharia e o Meio Ambiente; Os Ecossistemas. A crise energética. Fontes alternativas de energia. A sustentabilidade do meio ambiente. 2 - O MEIO AMBIENTE AQUÁTICO: Composição e Propriedades; Necessidade e Utilização; Requisitos de Qualidade; Poluição.3 - O MEIO AMBIENTE TERRESTRE: Composição e Propriedades; Necessidades e Utilização; Requisitos de Qualidade; Poluição.4 - O MEIO AMBIENTE ATMOSFÉRICO: Composição e Propriedades; Requisitos de Qualidade; Poluição."

# 6) Row 19 ("Método:") previously (wrongly) held the teacher's name - fix to
#    the evaluation method text.
$ws.Cells.Item(19, 2).Value = "Duas Provas  P1  1º bimestre e P2  2º bimestre"
$ws.Cells.Item(19, 3).Value = "Duas Provas  P1  1º bimestre e P2  2º bimestre"

# 7) Row 20 ("Critério:") previously held the evaluation-method text one row
#    early - fix to the grading criterion formula.
$ws.Cells.Item(20, 2).Value = "MF = (P1+ P2)/2"
$ws.Cells.Item(20, 3).Value = "MF = (P1+ P2)/2"

# 8) Row 21 ("Norma de recuperação:") previously held the criterion formula
#    one row early - fix to the recovery-norm text.
$ws.Cells.Item(21, 2).Value = "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação"
$ws.Cells.Item(21, 3).Value = "NF = (MF + PR)/ 2 , onde PR é uma prova de recuperação"

# 9) Row 22 ("Bibliografia:") previously held the recovery-norm text one row
#    early - fix to the actual bibliography (multi-line).
$bibliografia = "Braga, B.P.F., M.T.,Conejo, J.G., Porto, M.F., Veras M.S., Nucci, N., Juliano, N. e Eiger, S. Introdução à Engenharia Ambiental, Makron Books, São Paulo, 1998`nSperling, M.V. Princípios do Tratamento Biológico de Águas Residuárias. Desa-UFMG, Minas Gerais, 1996.`nBRAGA, B.et al. Introdução à Engenharia Ambiental. São Paulo: Prentice Hall, 2002, 305 p.`nVON SPERLING, M. Introdução à qualidade das águas e ao tratamento de esgotos. 2. ed. Belo Horizonte: UFMG, 1996."
$ws.Cells.Item(22, 2).Value = $bibliografia
$ws.Cells.Item(22, 3).Value = $bibliografia
